# ScenarioGuide.xlsx - "model set up with parallel processing"
#
# This script updates the four scenario-guide sheets:
#   - Para_Demo_Global (sheet1): a couple of parameter rows switch to
#     "Off" and have their step counts bumped.
#   - Coef_Demo_Global (sheet2): three coefficient rows get their
#     Min/Max collapsed onto the Default value (single-value run).
#   - Var_Output (sheet3): two new output columns (HistTable /
#     CompareSwitch) are added, driving the history-series comparisons.
#   - Test_Summary (sheet4): the run-time estimate is reworked to take
#     parallel CPUs into account (minutes/run, # CPUs, hours).

$wb = $excel.ActiveWorkbook

$wsPara = $wb.Worksheets.Item("Para_Demo_Global")
$wsCoef = $wb.Worksheets.Item("Coef_Demo_Global")
$wsVar  = $wb.Worksheets.Item("Var_Output")
$wsSum  = $wb.Worksheets.Item("Test_Summary")

# ---------------------------------------------------------------------
# Para_Demo_Global
# ---------------------------------------------------------------------
$wsPara.Range("E2").Value = "Off"
$wsPara.Range("I3").Value = 20
$wsPara.Range("H4").Value = 4
$wsPara.Range("I4").Value = 20
$wsPara.Range("E5").Value = "Off"

# ---------------------------------------------------------------------
# Coef_Demo_Global
# ---------------------------------------------------------------------
$wsCoef.Columns.Item(1).ColumnWidth = 11.333333333333334

$wsCoef.Range("I2").Value = 0.0109
$wsCoef.Range("J2").Value = 1

$wsCoef.Range("H3").Value = -0.8327
$wsCoef.Range("I3").Value = -0.8327
$wsCoef.Range("J3").Value = 1

$wsCoef.Range("H4").Value = -0.00948
$wsCoef.Range("I4").Value = -0.00948
$wsCoef.Range("J4").Value = 1

# ---------------------------------------------------------------------
# Var_Output / Test_Summary new columns & rows
#
# (String cells are entered in the order the author actually typed
# them so shared-string indices line up with the source workbook.)
# ---------------------------------------------------------------------
$wsVar.Range("D1").Value = "HistTable"
$wsSum.Range("A3").Value = "Total Number of CPUs (# of Model Runs in Parallel)"
$wsVar.Range("D3").Value = "SeriesPopulation"
$wsVar.Range("D2").Value = "SeriesTFR"
$wsVar.Range("D4").Value = "SeriesLifeExpectIHMEBothSexesHistOnly"
$wsVar.Range("D5").Value = "SeriesForecastBirthsMedUNPD"
$wsVar.Range("D6").Value = "SeriesDeathsper1000IHMEForecasts"
$wsVar.Range("E1").Value = "CompareSwitch"
$wsSum.Range("A5").Value = "Total Estimated Time (Hours)"
$wsSum.Range("A2").Value = "Mins per Model Run"

$wsVar.Range("E2").Value = "On"
$wsVar.Range("E3").Value = "On"
$wsVar.Range("E4").Value = "On"
$wsVar.Range("E5").Value = "On"
$wsVar.Range("E6").Value = "Off"

$wsVar.Columns.Item(4).ColumnWidth = 34.833333333333336
$wsVar.Columns.Item(5).ColumnWidth = 22.5

$wsVar.AutoFilterMode = $false
$wsVar.Range("A1:E6").AutoFilter()

$wsVar.Range("C2:C6,E2:E7").Validation.Delete()
$wsVar.Range("C2:C6,E2:E7").Validation.Add(3, 1, 1, "On,Off")

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Var_Output!_FilterDatabase") {
        $n.RefersTo = "=Var_Output!`$A`$1:`$E`$6"
    }
}

# ---------------------------------------------------------------------
# Test_Summary values / formulas / number formats
# ---------------------------------------------------------------------
$wsSum.Columns.Item(1).ColumnWidth = 42.666666666666664

$wsSum.Range("B4").Formula = "=(B1*(B2*60)/(60*60*24))/B3"
$wsSum.Range("B5").Formula = "=B4*24"

$wsSum.Range("B2").Value = 2.1
$wsSum.Range("B3").Value = 4

$wsSum.Range("B1").NumberFormat = "0"
$wsSum.Range("B3").NumberFormat = "0"
$wsSum.Range("B2").NumberFormat = "0.0"
$wsSum.Range("B4").NumberFormat = "0.0"
$wsSum.Range("B5").NumberFormat = "0.0"

$wsSum.Columns.Item(2).ColumnWidth = 8

# ---------------------------------------------------------------------
# Selections / active sheet (matches the author's last on-screen state)
# ---------------------------------------------------------------------
$wsCoef.Range("D13").Select()
$wsVar.Range("D17").Select()
$wsSum.Range("D11").Select()

$wsPara.Activate()
$wsPara.Range("E21").Select()

Write-Output "done"
